# ft && fix: removed not valid cuts && lp gap
#
# The match-schedule generator produced a corrected (valid) set of
# results for the "Resultados" sheet (some Local/Visita pairings were
# swapped and the score strings changed), which in turn changes the
# aggregated standings ("Puntos" / "Localias faltantes") on the
# "Equipos" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Equipos": update Puntos (D) / Localias faltantes (E) columns
# ---------------------------------------------------------------------
$wsEquipos = $wb.Worksheets.Item("Equipos")

$wsEquipos.Range("D2").Value = 4   # Team A points: 6 -> 4

$wsEquipos.Range("D3").Value = 6   # Team B points: 7 -> 6
$wsEquipos.Range("E3").Value = 2   # Team B localias faltantes: 1 -> 2

$wsEquipos.Range("D4").Value = 1   # Team C points: 4 -> 1
$wsEquipos.Range("E4").Value = 1   # Team C localias faltantes: 2 -> 1

$wsEquipos.Range("D5").Value = 6   # Team D points: 0 -> 6

# ---------------------------------------------------------------------
# Sheet "Resultados": update Local (C) / Visita (D) / Resultado (E)
# ---------------------------------------------------------------------
$wsResultados = $wb.Worksheets.Item("Resultados")

# Jornada 6
$wsResultados.Range("E3").Value = "8:4"   # was 2:3

$wsResultados.Range("E4").Value = "5:3"   # was 1:5

# Jornada 5
$wsResultados.Range("E6").Value = "2:2"   # was 6:4

$wsResultados.Range("C7").Value = "B"     # was D (swapped with D7)
$wsResultados.Range("D7").Value = "D"     # was B
$wsResultados.Range("E7").Value = "1:2"   # was 0:2

# Jornada 4
$wsResultados.Range("E9").Value = "3:4"   # was 1:4

$wsResultados.Range("C10").Value = "D"    # was C (swapped with D10)
$wsResultados.Range("D10").Value = "C"    # was D
$wsResultados.Range("E10").Value = "2:3"  # was 3:2

# Jornada 3
$wsResultados.Range("E12").Value = "4:3"  # was 0:6

$wsResultados.Range("C13").Value = "C"    # was B (swapped with D13)
$wsResultados.Range("D13").Value = "B"    # was C
$wsResultados.Range("E13").Value = "1:6"  # was 2:2

# Jornada 2
$wsResultados.Range("E15").Value = "2:2"  # was 5:4

$wsResultados.Range("E16").Value = "3:1"  # was 1:0

# Jornada 1
$wsResultados.Range("E18").Value = "5:0"  # was 3:4

$wsResultados.Range("E19").Value = "1:4"  # was 3:1
